# Se agrega funcionalidad para Menú
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("03-05-2022 08:44", "hola"),
    @("03-05-2022 08:45", "1"),
    @("03-05-2022 08:50", "hola"),
    @("03-05-2022 08:55", "hola"),
    @("03-05-2022 08:57", "hola"),
    @("03-05-2022 08:59", "hola"),
    @("03-05-2022 08:59", "hola"),
    @("03-05-2022 08:59", "quiero_info"),
    @("03-05-2022 09:03", "hola"),
    @("03-05-2022 09:03", "hola"),
    @("03-05-2022 09:04", "hola"),
    @("03-05-2022 09:04", "hola"),
    @("03-05-2022 09:04", "hola"),
    @("03-05-2022 09:04", "1"),
    @("03-05-2022 09:22", "hola"),
    @("03-05-2022 09:22", "1"),
    @("03-05-2022 09:23", "Hola"),
    @("03-05-2022 09:23", "Hola"),
    @("03-05-2022 09:23", "Hola"),
    @("03-05-2022 09:24", "Hola")
)

# "1" looks numeric, so a plain .Value assignment would store it as a
# number (losing the shared-string entry the diff expects, and it must be
# inserted into the shared-string table at the same position a literal
# text entry would take). Write it as genuine text without touching any
# cell's NumberFormat/style by building it via a TEXT() formula in a
# scratch cell, copying it, and paste-special'ing (values only) into the
# target cell - in the same left-to-right, top-to-bottom order as the rest
# of the data so new shared strings land in diff order.
$startRow = 11
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]

    if ($data[$i][1] -eq "1") {
        $ws.Range("Z1").Formula = "=TEXT(1,""0"")"
        $ws.Range("Z1").Copy()
        $ws.Cells.Item($row, 2).PasteSpecial(-4163)
    } else {
        $ws.Cells.Item($row, 2).Value = $data[$i][1]
    }
}
$ws.Range("Z1").Value = ""
